# Weekly fruit/vegetable price update: a new week's record is inserted
# at row 11 (pushing the existing rows 11-25 down to 12-26) and the new
# row is populated with the latest reading for Albahaca.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 11, shifting rows 11-25 down to 12-26.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with this week's data.
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C11").Value = "Arica y Parinacota"
$ws.Range("D11").Value = 44536
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 100112052
$ws.Range("G11").Value = "Albahaca"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 250
$ws.Range("K11").Value = 900
$ws.Range("L11").Value = 1000
$ws.Range("M11").Value = 950
$ws.Range("N11").Value = "$/paquete"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 950
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = "Hortaliza"
